$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "test"
$ws.Range("C3").Value = 612
$ws.Range("C5").Value = 69
$ws.Range("C7").Value = 25
